$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5706239938735962
$ws.Range("B1").Value = 2.78913950920105
$ws.Range("C1").Value = 3.174113750457764
$ws.Range("D1").Value = 3.79065728187561
$ws.Range("E1").Value = 1.121853590011597
